$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Лаба №3 (D) score corrected from -1 to 4
$ws.Range("D4").Value = 4

# Row 5: Лаба №5 (G) score entered as 0 (was blank)
$ws.Range("G5").Value = 0

# Row 7: Лаба №4 (E) score entered as -1 (was blank)
$ws.Range("E7").Value = -1

# Row 10: Лаба №4 (E) score entered as 4 (was blank)
$ws.Range("E10").Value = 4

# Row 14: Лаба №4 (E) score entered as -1 (was blank)
$ws.Range("E14").Value = -1

# Row 15: Лаба №2 (C) score corrected from 3 to 4 (keep the original cell
# formatting, which carries a quote-prefix style, by re-applying the
# format from a cell that shares the same style after updating the value)
$c15 = $ws.Range("C15")
$c15.Value = 4
$ws.Range("C6").Copy()
$c15.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 15: note changed to "переписаны верно все номера"
$ws.Range("M15").Value = "переписаны верно все номера"

# Row 17: Лаба №4 (E) score entered as -1 (was blank)
$ws.Range("E17").Value = -1

# Row 21: Лаба №4 (E) score entered as 4 (was blank)
$ws.Range("E21").Value = 4

# Row 22: КР-2 (F) score corrected from 4 to 3, and remark cleared
$ws.Range("F22").Value = 3
$ws.Range("N22").Value = ""

# Row 23: КР-2 (F) score corrected from 4 to 2, and remark cleared
$ws.Range("F23").Value = 2
$ws.Range("N23").Value = ""
